$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Save" column
$ws.Range("H1").Value = "Save"
$ws.Range("H1").Style = $ws.Range("G1").Style

# Save values per row (2..26), taken from source data
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
